$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "This sample is compatible with the Windows 10 Creators Update
#    SDK (15063)" -> "... Windows 10 Fall Creators Update SDK (16299)"
#    and move the hidden "_GoBack" bookmark to sit right after the new
#    text (this is what Word does automatically after the most recent
#    edit location).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Windows 10 Creators Update SDK (15063)", $false, $false, $false, $false, $false, $true, 1, $false, "Windows 10 Fall Creators Update SDK (16299)", 2)
$endPos = $rng.End

# Work around a collapsed-range-at-paragraph-end placement quirk: insert
# a throwaway character right after the new text, wrap the (non-empty)
# bookmark range around it, then delete the throwaway character again.
# The bookmark collapses back down to a zero-width bookmark that sits
# exactly after the inserted text.
$tmp = $d.Range($endPos, $endPos)
$tmp.InsertAfter("X")
$bmRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$delRange = $d.Range($endPos, $endPos + 1)
$delRange.Text = ""

# ------------------------------------------------------------------
# 2) Merge "<space>" + "See the " runs into a single "<space>See the "
#    run. (Round-trip through a scratch value first so the engine
#    actually records the run rewrite instead of treating an in-place
#    identical-text write as a no-op.)
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("outputs 8bit values using Rec.709 color primaries.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterPrimaries = $rng2.End
$mergeRange = $d.Range($afterPrimaries, $afterPrimaries + 9)
$mergeRange.Text = "XXSEEXX"
$mergeRange2 = $d.Range($afterPrimaries, $afterPrimaries + 7)
$mergeRange2.Text = " See the "

$word.ActiveDocument
